# "allow adding teams to groups in excel"
#
# - groups sheet: add a new row for the "auditors" team (marked "yes")
# - haruba sheet: add a new row granting the "adimian" group read/write
#   permissions ("Y")
# - update the selections left behind on the groups/teams/haruba sheets
#   so the workbook reopens focused on the newly added rows

$wb = $excel.ActiveWorkbook

# --- groups: register the "auditors" team -------------------------------
$groups = $wb.Worksheets.Item("groups")
$groups.Range("A4").Value = "auditors"
$groups.Range("B4").Value = "yes"

# --- haruba: grant the "adimian" group permissions -----------------------
$haruba = $wb.Worksheets.Item("haruba")
$haruba.Range("A5").Value = "adimian"
$haruba.Range("B5").Value = "Y"
$haruba.Range("C5").Value = "Y"

# --- teams: leave data untouched, just move the selection ---------------
$teams = $wb.Worksheets.Item("teams")
$teams.Select()
$teams.Range("B1:C1").Select()

# --- haruba: move selection onto the freshly added row -------------------
$haruba.Select()
$haruba.Range("A5:C6").Select()

# --- groups: move selection onto the freshly added row, and make sure
#     "groups" remains the tab that is active/selected when the workbook
#     is reopened (matches the original workbook state) ------------------
$groups.Select()
$groups.Range("A4").Select()
